# Generate Report for Handback
#
# The two files (44533c18... and 05ba5c4b...) have now come back from
# localization ("Handed back: in sync with en-US") - this script updates
# the Overview sheet and the per-language (zh-cn / de-de) detail sheets
# with the handback status/details, matching what a "generate handback
# report" run would produce.
#
# NOTE: this engine's PowerShell function calls do not bind named
# (-param value) arguments reliably, so every helper below uses purely
# positional parameters.

$wb = $excel.ActiveWorkbook

$file05 = "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md"
$file44 = "44533c18-eeca-4e90-af39-4b8a6044c26d.md"
$file50 = "50f396db-a00a-4d9f-bc95-ae77014c1455.md"
$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 now reports on 05ba5c4b...md (was 44533c18...md)
$ov.Range("A2").Value2 = $file05
$ov.Range("B2").Value2 = $statusHandedBack
$ov.Range("C2").Value2 = $statusHandedBack
$ov.Range("D2").Value2 = "2016-03-24 02:37:32"

# Row 3 now reports on 44533c18...md (was 05ba5c4b...md)
$ov.Range("A3").Value2 = $file44
$ov.Range("B3").Value2 = $statusHandedBack
$ov.Range("C3").Value2 = $statusHandedBack
$ov.Range("D3").Value2 = "2016-03-24 02:35:52"

# Row 4 (50f396db...md) moves from "Ready for handoff" to handed back
$ov.Range("A4").Value2 = $file50
$ov.Range("B4").Value2 = $statusHandedBack
$ov.Range("C4").Value2 = $statusHandedBack
$ov.Range("D4").Value2 = "2016-03-24 02:37:32"

# ---------------------------------------------------------------------
# Helper: add a hyperlink without disturbing cell style more than the
# engine already does, using purely positional args.
# ---------------------------------------------------------------------
function Add-FileLink {
    param($ws, $addr, $url, $display)
    $ws.Hyperlinks.Add($ws.Range($addr), $url, $null, $null, $display) | Out-Null
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$xlf05zh = "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf"
$xlf44zh = "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.zh-cn.xlf"

# --- Row 2: now the 05ba5c4b file (previously held 44533c18's data) ---
$zh.Range("A2").Value2 = $file05
$zh.Range("C2").Value2 = $statusHandedBack
$zh.Range("D2").Value2 = $xlf05zh
$zh.Range("E2").Value2 = "2016-03-24 02:37:23"
$zh.Range("F2").Value2 = $file05
$zh.Range("G2").Value2 = $xlf05zh
$zh.Range("H2").Value2 = "2016-03-24 02:38:20"

# --- Row 3: now the 44533c18 file (previously held 05ba5c4b's data) ---
$zh.Range("A3").Value2 = $file44
$zh.Range("C3").Value2 = $statusHandedBack
$zh.Range("D3").Value2 = $xlf44zh
$zh.Range("E3").Value2 = "2016-03-24 02:35:43"

# F3/G3 are brand-new cells (row previously had no Target/Handback file)
$zh.Range("F3").Value2 = $file44
$zh.Range("G3").Value2 = $xlf44zh
$zh.Range("H3").Value2 = "2016-03-24 02:36:23"

Add-FileLink $zh "F3" "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/26313cd274d994beb989eb4dd3ea733ae8b374b2/e2e/$file44" $file44
Add-FileLink $zh "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bd93f69531d589e46145e57657046d22b825ed41/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$xlf44zh" $xlf44zh

# --- Row 4: 50f396db file, still maps to the 05ba5c4b target/handback pair ---
$zh.Range("C4").Value2 = $statusHandedBack

$zh.Range("F4").Value2 = $file05
$zh.Range("G4").Value2 = $xlf05zh
$zh.Range("H4").Value2 = "2016-03-24 02:38:20"

Add-FileLink $zh "F4" "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/26313cd274d994beb989eb4dd3ea733ae8b374b3/e2e/$file05" $file05
Add-FileLink $zh "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bd93f69531d589e46145e57657046d22b825ed42/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/$xlf05zh" $xlf05zh

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$xlf05de = "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.de-de.xlf"
$xlf44de = "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.de-de.xlf"

# --- Row 2: now the 05ba5c4b file (previously held 44533c18's data) ---
$de.Range("A2").Value2 = $file05
$de.Range("C2").Value2 = $statusHandedBack
$de.Range("D2").Value2 = $xlf05de
$de.Range("E2").Value2 = "2016-03-24 02:37:32"
$de.Range("F2").Value2 = $file05
$de.Range("G2").Value2 = $xlf05de
$de.Range("H2").Value2 = "2016-03-24 02:38:36"

# --- Row 3: now the 44533c18 file (previously held 05ba5c4b's data) ---
$de.Range("A3").Value2 = $file44
$de.Range("C3").Value2 = $statusHandedBack
$de.Range("D3").Value2 = $xlf44de
$de.Range("E3").Value2 = "2016-03-24 02:35:52"

# F3/G3 are brand-new cells
$de.Range("F3").Value2 = $file44
$de.Range("G3").Value2 = $xlf44de
$de.Range("H3").Value2 = "2016-03-24 02:36:37"

Add-FileLink $de "F3" "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/cf51fb544b900cdc839933894f4487bcbf8ec4b6/e2e/$file44" $file44
Add-FileLink $de "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b0a11c532b16eb9ed54a221d97e9bf6c6d9c45b/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$xlf44de" $xlf44de

# --- Row 4: 50f396db file, still maps to the 05ba5c4b target/handback pair ---
$de.Range("C4").Value2 = $statusHandedBack

$de.Range("F4").Value2 = $file05
$de.Range("G4").Value2 = $xlf05de
$de.Range("H4").Value2 = "2016-03-24 02:38:36"

Add-FileLink $de "F4" "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/cf51fb544b900cdc839933894f4487bcbf8ec4b7/e2e/$file05" $file05
Add-FileLink $de "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b0a11c532b16eb9ed54a221d97e9bf6c6d9c45c/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/$xlf05de" $xlf05de

Write-Host "Handback report generated."
